$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
    "60-9="
    "80+13="
    "59-28="
    "78-33="
    "37+49="
    "70+6="
    "2+37="
    "72-37="
    "94-23="
    "10-5="
    "75-49="
    "55-51="
    "42-8="
    "79-1="
    "14+31="
    "84-24="
    "22+30="
    "35+14="
    "24+26="
    "81-3="
    "22+74="
    "15+17="
    "90-61="
    "6+7="
    "59+6="
    "5+29="
    "71-42="
    "0+22="
    "96-80="
    "9+72="
    "78+15="
    "90-32="
    "62+15="
    "33+20="
    "94-3="
    "59-4="
    "85-31="
    "23-13="
    "16+68="
    "26-6="
    "3+62="
    "57-34="
    "66-58="
    "69+11="
    "57+37="
    "50+1="
    "3+54="
    "9+29="
    "65-33="
    "81-10="
    "77-58="
    "62-20="
    "45-9="
    "71-10="
    "0+76="
    "3+49="
    "11+20="
    "6+85="
    "9+43="
    "96-3="
    "97-42="
    "21+18="
    "54+12="
    "33+41="
    "21+71="
    "0+77="
    "36+48="
    "11+5="
    "5+20="
    "1+53="
    "90-52="
    "9+67="
    "69-61="
    "36+46="
    "63-62="
    "3+49="
    "5+63="
    "29-13="
    "20+11="
    "89-24="
    "6+19="
    "92-5="
    "18-10="
    "24+35="
    "96-23="
    "13+43="
    "5+46="
    "51+31="
    "10+28="
    "17+78="
    "0+20="
    "74+10="
    "58+12="
    "87-84="
    "35+48="
    "31+40="
    "19-8="
    "4+38="
    "70-53="
    "28+12="
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count

if (($rows * $cols) -ne $newValues.Count) {
    throw "Cell count mismatch: table has $($rows * $cols) cells but $($newValues.Count) replacement values were supplied."
}

$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}

Write-Output "Updated $idx cells"